$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lesson 22 ("File I/O & Date-Time API") YouTube link: the Skype meeting
# link is replaced by the uploaded YouTube recording, with a real
# hyperlink (like the other lesson rows already have).
$ws.Hyperlinks.Add($ws.Range("F25"), "https://youtu.be/CPQoVMBvaPo", "", "", "https://youtu.be/CPQoVMBvaPo") | Out-Null

# Row 25 ("File I/O & Date-Time API") height tweak (14.9 -> 14.2), matching
# the other normal-height rows now that the row is "finished".
$ws.Rows.Item(25).RowHeight = 14.2

# Lesson name text tweaks for lessons 23 & 24.
$ws.Range("C26").Value = "Java 8 (Stream API)"
$ws.Range("C27").Value = "Java 8 continued, Section Project #1"

# Move the saved cursor/selection down to reflect the newly filled-in row.
$ws.Range("F31").Select() | Out-Null
